$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reverse the RT recursive / matching-layers loss series (flat + wedge) ---
# The existing 17-row series (rows 2..18) shifts down by one row (rows 3..19),
# a new leading data point is written into row 2, and two new trailing data
# points are appended as rows 20 and 21.

$ws.Cells.Item(2, 1).Value = -581.0752275067952
$ws.Cells.Item(2, 2).Value = -571.0445665845083

$ws.Cells.Item(3, 1).Value = -409.8017416462444
$ws.Cells.Item(3, 2).Value = -529.4950294068274

$ws.Cells.Item(4, 1).Value = -299.9771098537376
$ws.Cells.Item(4, 2).Value = -503.5727698408327

$ws.Cells.Item(5, 1).Value = -214.8498429659494
$ws.Cells.Item(5, 2).Value = -483.812940884937

$ws.Cells.Item(6, 1).Value = -142.8404221329019
$ws.Cells.Item(6, 2).Value = -467.2775617162167

$ws.Cells.Item(7, 1).Value = -78.8636054427202
$ws.Cells.Item(7, 2).Value = -452.6916457557847

$ws.Cells.Item(8, 1).Value = -20.24911567568273
$ws.Cells.Item(8, 2).Value = -439.3924802853804

$ws.Cells.Item(9, 1).Value = 34.57059100326809
$ws.Cells.Item(9, 2).Value = -426.9944869777443

$ws.Cells.Item(10, 1).Value = 86.58885582996362
$ws.Cells.Item(10, 2).Value = -415.2552718095986

$ws.Cells.Item(11, 1).Value = 136.4705680601647
$ws.Cells.Item(11, 2).Value = -404.0137595849461

$ws.Cells.Item(12, 1).Value = 184.6797209652547
$ws.Cells.Item(12, 2).Value = -393.1582975876381

$ws.Cells.Item(13, 1).Value = 231.5509384404663
$ws.Cells.Item(13, 2).Value = -382.6088781991294

$ws.Cells.Item(14, 1).Value = 277.3317785924383
$ws.Cells.Item(14, 2).Value = -372.3066757350266

$ws.Cells.Item(15, 1).Value = 322.209608050206
$ws.Cells.Item(15, 2).Value = -362.2074417896602

$ws.Cells.Item(16, 1).Value = 366.3287651220534
$ws.Cells.Item(16, 2).Value = -352.2772999699382

$ws.Cells.Item(17, 1).Value = 409.8022275632035
$ws.Cells.Item(17, 2).Value = -342.4899014874767

$ws.Cells.Item(18, 1).Value = 452.7196912703856
$ws.Cells.Item(18, 2).Value = -332.8244634297997

$ws.Cells.Item(19, 1).Value = 495.1532508711954
$ws.Cells.Item(19, 2).Value = -323.2643935625607

$ws.Cells.Item(20, 1).Value = 537.16149053607
$ws.Cells.Item(20, 2).Value = -313.7963030687208

$ws.Cells.Item(21, 1).Value = 578.7924866790967
$ws.Cells.Item(21, 2).Value = -305.9288606414241

# Column A carries a bold/bordered/centered style on every data row. Clone that
# formatting (via copy/paste-special of formats only, so no duplicate style
# entries get created) onto the brand-new column-A cells (row 2, 19, 20, 21 did
# not previously exist in the sheet, so a plain value write leaves them
# unstyled).
$ws.Range("A18").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
